# Update "想去人数" (want-to-go count) values across sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 285
$ws1.Range("F4").Value = 1204
$ws1.Range("F5").Value = 334
$ws1.Range("F6").Value = 122
$ws1.Range("F7").Value = 2985
$ws1.Range("F9").Value = 724
$ws1.Range("F10").Value = 485
$ws1.Range("F12").Value = 208
$ws1.Range("F13").Value = 721
$ws1.Range("F15").Value = 146
$ws1.Range("F16").Value = 1958
$ws1.Range("F17").Value = 316
$ws1.Range("F19").Value = 7

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 34
$ws2.Range("F9").Value = 123
$ws2.Range("F12").Value = 68

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6381
$ws3.Range("F3").Value = 805
$ws3.Range("F4").Value = 2053

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 0
$ws4.Range("F3").Value = 805
$ws4.Range("F4").Value = 2053
$ws4.Range("F9").Value = 34
$ws4.Range("F11").Value = 285
$ws4.Range("F12").Value = 1204
$ws4.Range("F13").Value = 334
$ws4.Range("F17").Value = 122
$ws4.Range("F18").Value = 2985
$ws4.Range("F19").Value = 123
$ws4.Range("F23").Value = 68
$ws4.Range("F24").Value = 724
$ws4.Range("F25").Value = 485
$ws4.Range("F26").Value = 319
$ws4.Range("F28").Value = 208
$ws4.Range("F29").Value = 721
$ws4.Range("F31").Value = 146
$ws4.Range("F33").Value = 1959
$ws4.Range("F34").Value = 316
$ws4.Range("F38").Value = 7
